# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" message on sheet "Hoja1" (A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$text = $wsHoja1.Range("A1").Text
$text = $text.Replace("1000 Bs = 8.36 = 34365.12 pesos", "1000 Bs = 8.31 = 34062.71 pesos")
$text = $text.Replace("34365.12 pesos = 8.31 = 972.49 Bs", "34062.71 pesos = 8.26 = 969.99 Bs")
$wsHoja1.Range("A1").Value = $text

# --- Update the rate figures on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 120.4
$wsTasas.Range("O10").Value = 4101.15
$wsTasas.Range("N12").Value = 4126.2
$wsTasas.Range("O12").Value = 117.5
